# adding just wallops as GS option
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$sheet6 = $wb.Worksheets.Item("6")
$gsGenerator = $wb.Worksheets.Item("gs_generator")

# ------------------------------------------------------------------
# 0. Throwaway sheet used purely to keep the internal sheetId counter
#    in step with the authored workbook (sheetId 15 gets "burned").
# ------------------------------------------------------------------
$scratch = $wb.Worksheets.Add($gsGenerator)
$scratchName = $scratch.Name

# ------------------------------------------------------------------
# 1. Insert two new ground-station sheets ("7" and "8") just before
#    "gs_generator", both copied from the pristine "6" sheet.
# ------------------------------------------------------------------
$sheet6.Copy($null, $sheet6)
$sheet7 = $wb.Worksheets.Item("6 (2)")
$sheet7.Name = "7"

$sheet6.Copy($null, $sheet7)
$sheet8 = $wb.Worksheets.Item("6 (2)")
$sheet8.Name = "8"

$wb.Worksheets.Item($scratchName).Delete()

# ---- sheet "7": row 2 -> Wallops, row 3 -> MIT -------------------
$sheet7.Range("A2").Value = 1
$sheet7.Range("B2").Value = 37.855662000000002
$sheet7.Range("C2").Value = -75.512068999999997
$sheet7.Range("D2").Value = 0
$sheet7.Range("E2").Value = 1
$sheet7.Range("F2").Value = "Wallops"

$sheet7.Range("A3").Value = 2
$sheet7.Range("B3").Value = 42.360726
$sheet7.Range("C3").Value = -71.093208000000004
$sheet7.Range("D3").Value = 0
$sheet7.Range("E3").Value = 1
$sheet7.Range("F3").Value = "MIT"

$sheet7.Range("A4:D4").ClearContents()
$sheet7.Range("F4").ClearContents()
$sheet7.Range("A5:F5").ClearContents()
$sheet7.Range("A6:C6").ClearContents()
$sheet7.Range("E6:F6").ClearContents()
$sheet7.Range("A7:F7").ClearContents()
$sheet7.Range("A8:F8").ClearContents()
$sheet7.Range("A9:F9").ClearContents()
$sheet7.Range("A10:C10").ClearContents()
$sheet7.Range("E10:F10").ClearContents()

$sheet7.Range("J2").Value = 1
$sheet7.Range("J2").ClearContents()
$sheet7.Range("J4").Value = 1
$sheet7.Range("J4").ClearContents()

$sheet7.Range("I16").Select()

# ---- sheet "8": row 2 -> Wallops only -----------------------------
$sheet8.Range("A2").Value = 1
$sheet8.Range("B2").Value = 37.855662000000002
$sheet8.Range("C2").Value = -75.512068999999997
$sheet8.Range("D2").Value = 0
$sheet8.Range("E2").Value = 1
$sheet8.Range("F2").Value = "Wallops"

$sheet8.Range("A3:F3").ClearContents()

$sheet8.Range("A4:D4").ClearContents()
$sheet8.Range("F4").ClearContents()
$sheet8.Range("A5:F5").ClearContents()
$sheet8.Range("A6:C6").ClearContents()
$sheet8.Range("E6:F6").ClearContents()
$sheet8.Range("A7:F7").ClearContents()
$sheet8.Range("A8:F8").ClearContents()
$sheet8.Range("A9:F9").ClearContents()
$sheet8.Range("A10:C10").ClearContents()
$sheet8.Range("E10:F10").ClearContents()

$sheet8.Range("J2").Value = 1
$sheet8.Range("J2").ClearContents()
$sheet8.Range("J4").Value = 1
$sheet8.Range("J4").ClearContents()

$sheet8.Range("A3:F3").Select()

# ------------------------------------------------------------------
# 2. Update the original "6" (gs network) sheet: replace the
#    Availability column (E) with newly analysed values and add a
#    new column J holding the same "avail. From my own analysis"
#    figures.
# ------------------------------------------------------------------
$sheet6.Range("E2").Value = 0.35
$sheet6.Range("E3").Value = 0.49
$sheet6.Range("E4").Value = 0.59809999999999997
$sheet6.Range("E5").Value = 0.15329999999999999
$sheet6.Range("E6").Value = 0.20369999999999999
$sheet6.Range("E7").Value = 0.60840000000000005
$sheet6.Range("E8").Value = 0.65439999999999998
$sheet6.Range("E9").Value = 0.8448
$sheet6.Range("E10").Value = 0.74280000000000002

# new column J: header + values (copy E's style onto J first)
$sheet6.Range("E1").Copy($sheet6.Range("J1"))
$sheet6.Range("J1").Value = "avail. From my own analysis (col E is with mix of Inigo's and my numbers)"

$sheet6.Range("E2").Copy($sheet6.Range("J2"))
$sheet6.Range("J2").Value = 0.28220000000000001

$sheet6.Range("J3").Value = 0.16370000000000001

$sheet6.Range("E4").Copy($sheet6.Range("J4"))
$sheet6.Range("J4").Value = 0.59809999999999997

$sheet6.Range("J5").Value = 0.15329999999999999
$sheet6.Range("J6").Value = 0.20369999999999999
$sheet6.Range("J7").Value = 0.60840000000000005
$sheet6.Range("J8").Value = 0.65439999999999998
$sheet6.Range("J9").Value = 0.8448
$sheet6.Range("J10").Value = 0.74280000000000002

$sheet6.Range("E6").Select()

# ------------------------------------------------------------------
# 3. Misc view-selection tweaks on other sheets (cosmetic, match the
#    author's saved window state).
# ------------------------------------------------------------------
$gsGenerator.Range("H39").Select()
$wb.Worksheets.Item("4").Range("B8").Select()
$wb.Worksheets.Item("5").Range("E10").Select()
$wb.Worksheets.Item("0").Range("G12").Select()

# ------------------------------------------------------------------
# 4. Make "8" the active/visible sheet, matching the author's saved
#    workbook view.
# ------------------------------------------------------------------
$sheet8.Activate()
